$d = $word.ActiveDocument

$newParagraphs = @(
    "2/8/2022: Merged all branches and created a document to start tracking bugs and how to simulate them.",
    "2/9/2022: Created a functioning demo and merged everyone’s branches together.",
    "2/10/2022: Completed Sprint demo",
    "2/24/2022: Integrated the wiring scene and the bot customization scene.",
    "3/2/2022: Added keys to the mesh around the boe-bot so they can be referenced. Worked on the SRS and adding adits from Dr.Akbas.",
    "3/4/2022: Was able to get the wiring change to work with the left face of the bot-bot",
    "3/5/2022: Added wiring change to all faces of the boe-bot. Merged Daniel’s branch into my branch and fixed problems caused by merge. Added wiring change to the IR sensors. Added pins to all sensors used in the simulation.",
    "3/6/2022: Finished the SDS and fixed pins from flying around the scene. ",
    "3/9/2022: Working on fixing wiring scene bug. Merged Daniel’s branch. Created unity executables. Merged Arduino interface and added panning in the wiring scene.",
    "3/10/2022: Working on the wiring scene bugs.",
    "3/21/2022: Fixed the wiring scene bug with adding components late and switching scenes. Fixed bug where a user could attach a sensor in the wiring scene. "
)

foreach ($text in $newParagraphs) {
    $p = $d.Paragraphs.Add()
    $p.Range.Text = $text
}
